$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: worked 3h instead of 4h (running total recalculates automatically) ---
$ws.Range("B17").Value = 3

# --- New row 18: new entry for 2020-12-18 ---
# Copy formatting (number format / wrap style) from the row above so the
# new cells reuse the existing style entries instead of minting new ones.
$ws.Range("A17").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = 44183

$ws.Range("B18").Value = 4

$ws.Range("C18").Formula = "=B18+C17"

$ws.Range("D17").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "Interne Absprache, Planung weiteres vorghen`nQuellcodes auf neuen Redpitaya geschrieben und Funktion getestet"

# Match the real Excel autofit height (15pt/line * 3 lines) used elsewhere
# in this sheet for wrapped, multi-line entries.
$ws.Rows.Item(18).RowHeight = 45

# --- Update the current selection to reflect where the editor left off ---
$ws.Range("C21").Select() | Out-Null
